# Apply targeted updates to the "dSF" (F) column values on Sheet1.
# These correspond to a re-pull of source data / recalculated mean values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "F3"  = -3
    "F13" = -4
    "F14" = -5
    "F15" = -2
    "F20" = -8
    "F21" = -7
    "F22" = -2
    "F23" = -2
    "F24" = -2
    "F26" = -4
    "F31" = -14
    "F35" = -7
    "F36" = 3
    "F41" = -6
    "F51" = -5
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
